# aggiornamento fino a 21 marzo
# Append 4 new daily rows (230-233) to the end of the existing data table,
# continuing the date/contagi series already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append, in the same column order as the existing data:
# A = date serial, B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44304, 2, 15, 148.9572989076465),
    @(44305, 3, 15, 148.9572989076465),
    @(44306, 2, 15, 148.9572989076465),
    @(44307, 6, 21, 208.540218470705)
)

$lastRow = 229
$firstNewRow = $lastRow + 1
$lastNewRow = $lastRow + $newRows.Count

# Copy the date-column formatting (style) from the previous last row so the
# new date cells keep the same date number format / alignment / border.
$ws.Range("A$lastRow").Copy()
$ws.Range("A${firstNewRow}:A${lastNewRow}").PasteSpecial(-4122)

$r = $firstNewRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
